$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.538.36'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '2.939.15'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.54'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.65'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.30%  '
$ws.Range("E9").Value = '  +1.77%  '
$ws.Range("D10").Value = '2.937.48'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.442'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.89'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '3.476.02'
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.47'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '76.463.58'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = '2.935.13'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.53'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.74'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.62'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.27'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.87'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").Value = '3.082.48'
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.28'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.67'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000108'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.36'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.30%  '
$ws.Range("E32").Value = '  -3.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '499.60'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.36'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.16'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.392'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.61%  '
$ws.Range("E39").Value = '  +17.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.97'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.40%  '
$ws.Range("E41").Value = '  -4.37%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '179.55'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.92'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("E45").Value = '  -2.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.96'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("E47").Value = '  -4.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.593'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.87'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.663'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.73%  '
